# Auto update stock data: update the "Date_1" value (cell A2) on every
# worksheet from 2025/12/31 to 2025/10/24. The date is stored as plain
# text (not an Excel date serial), so force a text number format before
# assigning the new value to keep Excel from reinterpreting the string
# as a date.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A2")
    $cell.NumberFormat = "@"
    $cell.Value = "2025/10/24"
}
